{"js": "// Word JS API (Office.js) edit script.\n// Applies three changes to the letter:\n//   1. Update the letter date from \"September 19, 2025\" to \"September 21, 2025\".\n//   2. Split the single-line mailing address \"3503 Toomey Pl, Santa Clara CA 95051\"\n//      (the standalone paragraph under the addressee name, not the matching text\n//      that also lives inside the account-summary table) into two paragraphs:\n//      \"3503 Toomey Pl\" and \"Santa Clara, CA 95051\".\n//   3. Remove the empty \"No Spacing\" paragraph that immediately follows the\n//      \"...Board of Directors\" signature line.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Update the date line.\n// ---------------------------------------------------------------------\nconst dateResults = body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Split the mailing-address paragraph into street + city/state/zip.\n//    First figure out (read-only) which match index is the standalone\n//    paragraph rather than the one inside the account-summary table, then\n//    re-resolve that match with a fresh search right before mutating it\n//    (mutating earlier search hits can invalidate later ones in this\n//    runtime, so we avoid touching to-be-mutated proxies across a sync).\n// ---------------------------------------------------------------------\nconst addressText = \"3503 Toomey Pl, Santa Clara CA 95051\";\nconst addressResults = body.search(addressText, { matchCase: true });\naddressResults.load(\"items\");\nawait context.sync();\n\nlet addressIndex = -1;\nfor (let i = 0; i < addressResults.items.length; i++) {\n  const hitParagraphs = addressResults.items[i].paragraphs;\n  hitParagraphs.load(\"items\");\n  await context.sync();\n\n  const inTable = hitParagraphs.items[0].parentTableCellOrNullObject;\n  inTable.load(\"isNullObject\");\n  await context.sync();\n\n  if (inTable.isNullObject) {\n    addressIndex = i;\n    break;\n  }\n}\n\nif (addressIndex >= 0) {\n  const freshAddressResults = body.search(addressText, { matchCase: true });\n  freshAddressResults.load(\"items\");\n  await context.sync();\n\n  const hit = freshAddressResults.items[addressIndex];\n  const paragraph = hit.paragraphs.items[0];\n\n  // Replace the paragraph text with just the street address...\n  hit.insertText(\"3503 Toomey Pl\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // ...then add a new paragraph directly after it with the city/state/zip.\n  paragraph.insertParagraph(\"Santa Clara, CA 95051\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Delete the stray empty paragraph right after \"...Board of Directors\".\n// ---------------------------------------------------------------------\nconst boardResults = body.search(\"Board of Directors\", { matchCase: true });\nboardResults.load(\"items\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardParagraphs = boardResults.items[0].paragraphs;\n  boardParagraphs.load(\"items\");\n  await context.sync();\n\n  const boardParagraph = boardParagraphs.items[0];\n  const nextParagraph = boardParagraph.getNextOrNullObject();\n  nextParagraph.load(\"isNullObject,text,style\");\n  await context.sync();\n\n  if (!nextParagraph.isNullObject && nextParagraph.text === \"\" && nextParagraph.style === \"No Spacing\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM-interop edit script.\n# Applies three changes to the letter:\n#   1. Update the letter date from \"September 19, 2025\" to \"September 21, 2025\".\n#   2. Split the single-line mailing address \"3503 Toomey Pl, Santa Clara CA 95051\"\n#      (the standalone paragraph under the addressee name, not the matching text\n#      that also lives inside the account-summary table) into two paragraphs:\n#      \"3503 Toomey Pl\" and \"Santa Clara, CA 95051\".\n#   3. Remove the empty \"No Spacing\" paragraph that immediately follows the\n#      \"...Board of Directors\" signature line.\n\n$d = $word.ActiveDocument\n$wdWithInTable = 12\n$wdCollapseEnd = 0\n\n# ---------------------------------------------------------------------\n# 1) Update the date line.\n# ---------------------------------------------------------------------\n$dateRange = $d.Content\n$dateRange.Find.Execute(\"September 19, 2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"September 21, 2025\", 2)\n\n# ---------------------------------------------------------------------\n# 2) Split the mailing-address paragraph into street + city/state/zip.\n#    Walk every match of the address text and only touch the one that is\n#    not inside the account-summary table.\n# ---------------------------------------------------------------------\n$addressText = \"3503 Toomey Pl, Santa Clara CA 95051\"\n$searchRange = $d.Content\n$searchRange.Start = 0\n$searchRange.End = $d.Content.End\n\nwhile ($searchRange.Find.Execute($addressText)) {\n    $inTable = $searchRange.Information($wdWithInTable)\n    if (-not $inTable) {\n        $addrPara = $searchRange.Paragraphs(1)\n        $paraRange = $addrPara.Range\n        # Replace the visible text (leave the trailing paragraph mark alone).\n        $textRange = $d.Range($paraRange.Start, $paraRange.End - 1)\n        $textRange.Text = \"3503 Toomey Pl\"\n\n        # Insert a new paragraph right after it with the city/state/zip.\n        $addrPara2 = $d.Range($textRange.Start, $textRange.End).Paragraphs(1)\n        $addrPara2.Range.InsertParagraphAfter()\n        $addrPara2.Next().Range.Text = \"Santa Clara, CA 95051\"\n        break\n    }\n    $searchRange.Collapse($wdCollapseEnd)\n}\n\n# ---------------------------------------------------------------------\n# 3) Delete the stray empty paragraph right after \"...Board of Directors\".\n# ---------------------------------------------------------------------\n$boardRange = $d.Content\nif ($boardRange.Find.Execute(\"Board of Directors\")) {\n    $boardPara = $boardRange.Paragraphs(1)\n    $nextPara = $boardPara.Next()\n    if ($nextPara.Range.Text -eq [char]13 -and $nextPara.Style.NameLocal -eq \"No Spacing\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
